# Updated legacy GSC export files: drop the two oldest daily rows from the
# "Chart" sheet (2025-09-24 and 2025-09-25), shifting the remaining rows up.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows("2:3").Delete()
